$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '26.618.91'
$ws.Range('E2').Value = '  -2.19%  '
$ws.Range('D3').Value = '1.585.29'
$ws.Range('E3').Value = '  -2.81%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('E5').Value = '  -2.30%  '
$ws.Range('E6').Value = '  -2.21%  '
$ws.Range('E7').Value = '  -0.01%  '
$ws.Range('E8').Value = '  -2.52%  '
$ws.Range('E9').Value = '  -0.99%  '
$ws.Range('D10').Value = "'19.58"
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  -3.26%  '
$ws.Range('D11').Value = "'0.0833"
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  -1.67%  '
$ws.Range('D12').Value = '1.805.50'
$ws.Range('E12').Value = '  -2.84%  '
$ws.Range('D13').Value = '1.570.90'
$ws.Range('E13').Value = '  -3.60%  '
$ws.Range('D14').Value = "'4.05"
$ws.Range('D14').ClearFormats()
$ws.Range('E14').Value = '  -1.69%  '
$ws.Range('E15').Value = '  -2.88%  '
$ws.Range('D16').Value = "'64.50"
$ws.Range('D16').ClearFormats()
$ws.Range('D17').Value = '26.623.20'
$ws.Range('E17').Value = '  -2.04%  '
$ws.Range('E18').Value = '  -0.51%  '
$ws.Range('E19').Value = '  +0.03%  '
$ws.Range('D20').Value = "'207.73"
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  -3.92%  '
$ws.Range('E21').Value = '  -2.22%  '
$ws.Range('E22').Value = '  -3.12%  '
$ws.Range('D23').Value = "'2.36"
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  -5.02%  '
$ws.Range('D24').Value = "'8.90"
$ws.Range('D24').ClearFormats()
$ws.Range('E24').Value = '  -2.28%  '
$ws.Range('D25').Value = "'146.69"
$ws.Range('D25').ClearFormats()
$ws.Range('E25').Value = '  -0.85%  '
$ws.Range('E26').Value = '  -0.05%  '
$ws.Range('E27').Value = '  +1.64%  '
$ws.Range('E28').Value = '  -4.11%  '
$ws.Range('D29').Value = "'15.30"
$ws.Range('D29').ClearFormats()
$ws.Range('E31').Value = '  -1.98%  '
$ws.Range('E32').Value = '  -4.15%  '
$ws.Range('D33').Value = "'0.679"
$ws.Range('D33').ClearFormats()
$ws.Range('E33').Value = '  +25.25%  '
$ws.Range('D35').Value = '1.326.44'
$ws.Range('E35').Value = '  +0.80%  '
$ws.Range('D36').Value = "'1.51"
$ws.Range('D36').ClearFormats()
$ws.Range('E36').Value = '  -3.46%  '
$ws.Range('D37').Value = "'2.43"
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  -0.91%  '
$ws.Range('D39').Value = "'0.825"
$ws.Range('D39').ClearFormats()
$ws.Range('E39').Value = '  -3.04%  '
$ws.Range('E40').Value = '  +0.00%  '
$ws.Range('E41').Value = '  +3.05%  '
$ws.Range('E42').Value = '  -2.18%  '
$ws.Range('D43').Value = "'2.17"
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  -3.71%  '
$ws.Range('D44').Value = "'63.54"
$ws.Range('D44').ClearFormats()
$ws.Range('E44').Value = '  -0.12%  '
$ws.Range('D45').Value = '1.719.55'
$ws.Range('E45').Value = '  -2.64%  '
$ws.Range('D46').Value = "'89.64"
$ws.Range('D46').ClearFormats()
$ws.Range('E46').Value = '  -1.18%  '
$ws.Range('E47').Value = '  +1.06%  '
$ws.Range('E49').Value = '  +3.10%  '
$ws.Range('D51').Value = "'7.45"
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -1.11%  '
